# "Generate Report for Handoff"
#
# Refreshes the localization-status report: the "3a73c11f..." row moves
# from "handed back" to "ready for handoff" again (status text + the two
# handoff/handback timestamps on it get bumped), the zh-cn locale flag
# flips from "ht" to "mt", and a version-mismatch error is recorded for
# that row's handback file in both the zh-cn and de-de tabs. A couple of
# columns (Status-ish datetime column / Error Detail column) are also
# resized so the new content is readable.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------
# Overview sheet: zh-cn / de-de status + "Latest HO Xliff Generate Date"
# for the 3a73c11f-... row (row 2) and the 5c98e9ba-... row (row 3) both
# showed the same shared status/date text, so both rows move together.
# ---------------------------------------------------------------------
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-31 13:09:19"

$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-31 13:09:19"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("C3").Value = "Ready for handoff"

$zhcn.Range("E2").Value = "mt"
$zhcn.Range("E3").Value = "mt"

$zhcn.Range("H2").Value = "2016-08-31 13:09:02"
$zhcn.Range("H3").Value = "2016-08-31 13:09:02"

$zhcn.Range("P2").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0f18abe53d2f88cd3ed20620061cc82a66666e7c/e2e/3a73c11f-eca7-41bf-9da7-aa9e86668101.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/57aa5f56e4ecec0d882969e6c1c8b6505d316124/e2e/3a73c11f-eca7-41bf-9da7-aa9e86668101.md."

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("C3").Value = "Ready for handoff"

$dede.Range("E2").Value = "mt"
$dede.Range("E3").Value = "mt"

$dede.Range("H2").Value = "2016-08-31 13:09:19"
$dede.Range("H3").Value = "2016-08-31 13:09:19"

$dede.Range("P2").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0f18abe53d2f88cd3ed20620061cc82a66666e7c/e2e/3a73c11f-eca7-41bf-9da7-aa9e86668101.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/57aa5f56e4ecec0d882969e6c1c8b6505d316124/e2e/3a73c11f-eca7-41bf-9da7-aa9e86668101.md."

# ---------------------------------------------------------------------
# Column width tweaks that came along with this report refresh. The host
# app quantizes ColumnWidth to whole pixels and re-offsets it on save, so
# the assigned values are pre-compensated (target - 5/6) to land the
# saved sheet width on the closest representable value to the target
# (17.2159881591797 and 40 respectively).
# ---------------------------------------------------------------------
$overview.Columns.Item(5).ColumnWidth = 16.38265482584637
$overview.Columns.Item(6).ColumnWidth = 16.38265482584637

$zhcn.Columns.Item(3).ColumnWidth = 16.38265482584637
$zhcn.Columns.Item(16).ColumnWidth = 39.166666666666664

$dede.Columns.Item(3).ColumnWidth = 16.38265482584637
$dede.Columns.Item(16).ColumnWidth = 39.166666666666664
